$d = $word.ActiveDocument

$d.Content.Find.Execute("161×4=644", $true, $false, $false, $false, $false, $true, 1, $false, "241×4=964", 2)
$d.Content.Find.Execute("643×4=2572", $true, $false, $false, $false, $false, $true, 1, $false, "318×2=636", 2)
$d.Content.Find.Execute("658×7=4606", $true, $false, $false, $false, $false, $true, 1, $false, "225×5=1125", 2)
$d.Content.Find.Execute("570×7=3990", $true, $false, $false, $false, $false, $true, 1, $false, "451×2=902", 2)
$d.Content.Find.Execute("156×6=936", $true, $false, $false, $false, $false, $true, 1, $false, "822×6=4932", 2)
$d.Content.Find.Execute("779×9=7011", $true, $false, $false, $false, $false, $true, 1, $false, "113×2=226", 2)
$d.Content.Find.Execute("655×6=3930", $true, $false, $false, $false, $false, $true, 1, $false, "712×2=1424", 2)
$d.Content.Find.Execute("859×5=4295", $true, $false, $false, $false, $false, $true, 1, $false, "255×5=1275", 2)
$d.Content.Find.Execute("905×2=1810", $true, $false, $false, $false, $false, $true, 1, $false, "470×6=2820", 2)
$d.Content.Find.Execute("744×6=4464", $true, $false, $false, $false, $false, $true, 1, $false, "571×7=3997", 2)
$d.Content.Find.Execute("299×8=2392", $true, $false, $false, $false, $false, $true, 1, $false, "517×6=3102", 2)
$d.Content.Find.Execute("933×2=1866", $true, $false, $false, $false, $false, $true, 1, $false, "773×5=3865", 2)
$d.Content.Find.Execute("249×6=1494", $true, $false, $false, $false, $false, $true, 1, $false, "141×5=705", 2)
$d.Content.Find.Execute("665×5=3325", $true, $false, $false, $false, $false, $true, 1, $false, "745×4=2980", 2)
$d.Content.Find.Execute("440×7=3080", $true, $false, $false, $false, $false, $true, 1, $false, "390×5=1950", 2)
$d.Content.Find.Execute("515×7=3605", $true, $false, $false, $false, $false, $true, 1, $false, "307×3=921", 2)
$d.Content.Find.Execute("486×7=3402", $true, $false, $false, $false, $false, $true, 1, $false, "343×2=686", 2)
$d.Content.Find.Execute("975×2=1950", $true, $false, $false, $false, $false, $true, 1, $false, "954×2=1908", 2)
$d.Content.Find.Execute("545×5=2725", $true, $false, $false, $false, $false, $true, 1, $false, "465×5=2325", 2)
$d.Content.Find.Execute("233×6=1398", $true, $false, $false, $false, $false, $true, 1, $false, "883×7=6181", 2)
$d.Content.Find.Execute("528×4=2112", $true, $false, $false, $false, $false, $true, 1, $false, "643×4=2572", 2)
$d.Content.Find.Execute("776×2=1552", $true, $false, $false, $false, $false, $true, 1, $false, "991×6=5946", 2)
$d.Content.Find.Execute("938×4=3752", $true, $false, $false, $false, $false, $true, 1, $false, "829×7=5803", 2)
$d.Content.Find.Execute("122×8=976", $true, $false, $false, $false, $false, $true, 1, $false, "475×8=3800", 2)
$d.Content.Find.Execute("114×7=798", $true, $false, $false, $false, $false, $true, 1, $false, "828×8=6624", 2)
